# "added ppt slide for context"
#
# 1. Insert a brand-new slide at position 1 using the "Title and Content"
#    custom layout (same layout already used by the rest of the deck),
#    and fill it in with the situation / question set ("Set up" slide).
# 2. Move the (hidden) "To reduce risk, exercise caution around dangerous
#    breeds" slide from its old position (3) to the very end of the deck.

$p = $ppt.ActivePresentation

# Re-use the "Title and Content" layout already applied to slide 2 so the
# new slide gets a plain title + body placeholder (not the title-slide
# ctrTitle/subTitle pairing used by slide 1).
$titleAndContent = $p.Slides.Item(2).CustomLayout

$setup = $p.Slides.AddSlide(1, $titleAndContent)

# --- Title -----------------------------------------------------------
$setup.Shapes.Item(1).TextFrame.TextRange.Text = "Set up "

# --- Body / content placeholder --------------------------------------
$body = $setup.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$lines = @(
    "Situation: I work for the mayor’s office in the city of Louisville and they need to make a public statement in response to a number of recent dog bite incidents in the city.  ",
    "Questions under consideration",
    "Are the number of dog bites increasing, decreasing, or staying the same?",
    "And how does that compare to the number of dogs or rate of dog ownership in the US?",
    "Are you more likely to be bitten by a male dog or a female dog?",
    "Are dog bites concentrated in certain geographic areas?",
    "Which breed of dog are you most likely to be bitten by?",
    "I once lived in an apartment where certain dogs were not allowed.  Are they actually a greater risk for biting people?",
    "Are pit bulls more likely to bite people than other dog breeds?",
    "What are the chances you will get rabies from a dog bite?"
)

# Build the paragraphs up one at a time (rather than one big Text=
# assignment) so every run keeps its own run-level properties (lang
# etc.) instead of being collapsed into a single un-tagged run.
$tr.Text = $lines[0]
for ($i = 1; $i -lt $lines.Count; $i++) {
    [void]$tr.InsertAfter("`r" + $lines[$i])
}

# Paragraphs 3-10 (the "questions" bullets) are demoted one level and
# rendered smaller / italic.
for ($i = 3; $i -le $lines.Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.IndentLevel = 2
    $para.Font.Size = 14
    $para.Font.Italic = $true
}

# --- Move the "dangerous breeds" slide to the end of the deck --------
# After inserting the new slide at position 1 everything shifted down by
# one, so the slide that used to be 3rd ("To reduce risk, exercise
# caution around dangerous breeds") is now 4th.
$dangerousBreeds = $p.Slides.Item(4)
$dangerousBreeds.MoveTo($p.Slides.Count)
